$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Update the genBankAssembly accession for the outgroup row (row 9, col C)
#    from the old WGS accession prefix to the NCBI Assembly accession.
$ws.Range("C9").Value = "GCA_000703365.1"

# 2. Re-sort the data rows (A9:J18) ascending by biosample_acc (column A),
#    matching a user selecting the table and choosing Data > Sort.
$dataRange = $ws.Range("A9:J18")
$ws.Sort.SortFields.Clear()
$ws.Sort.SortFields.Add($ws.Range("A9:A18"))
$ws.Sort.SetRange($dataRange)
$ws.Sort.Header = -4142
$ws.Sort.Apply()

# 3. Widen column H slightly so the longer sha256sumAssembly values are legible.
$ws.Columns.Item(8).ColumnWidth = 21.59

# 4. Leave the cursor where the user last clicked after scrolling the sorted table.
$ws.Range("E25").Select()
